$d = $word.ActiveDocument

$d.Content.Find.Execute("Welcome to ParentText! ", $true, $false, $false, $false, $false, $true, 1, $false, "Welkom in ParentText! ", 2) | Out-Null
$d.Content.Find.Execute("Welcome to ParentText!", $true, $false, $false, $false, $false, $true, 1, $false, "Welkom by ParentText!", 2) | Out-Null
$d.Content.Find.Execute("Welcome to ParentText", $true, $false, $false, $false, $false, $true, 1, $false, "Welkom by ParentText", 2) | Out-Null
$d.Content.Find.Execute("ParentText is like having a supportive friend by your side, guiding you through the ups and downs of raising your child. ", $true, $false, $false, $false, $false, $true, 1, $false, "ParentText is soos om 'n ondersteunende vriend langs jou te hê, wat jou lei deur die op- en af drande van jou kind grootmaak. ", 2) | Out-Null
$d.Content.Find.Execute("Created by a team of experts from the University of Cape Town, the University of Oxford, Fort Hare University, Clowns Without Borders South Africa, and Parenting for Lifelong Health, ParentText has been tested worldwide to ensure it offers the best help possible. ", $true, $false, $false, $false, $false, $true, 1, $false, "Geskep deur 'n span kundiges van die Universiteit van Kaapstad, die Universiteit van Oxford, Fort Hare Universiteit, Clowns Without Borders Suid-Afrika, en Parenting for Lifelong Health, is ParentText wêreldwyd getoets om te verseker dat dit die beste hulp moontlik bied. ", 2) | Out-Null
$d.Content.Find.Execute("I am ______, your guide. Even though I might seem human, I am a robot created by Parenting for Lifelong Health and UNICEF to support you in your parenting journey. ", $true, $false, $false, $false, $false, $true, 1, $false, "Ek is ______, jou gids. Al lyk ek dalk menslik, is ek 'n robot, geskep deur Parenting for Lifelong Health en UNICEF om jou in jou ouerskapreis te ondersteun. ", 2) | Out-Null
$d.Content.Find.Execute("Let us see how ParentText works. ", $true, $false, $false, $false, $false, $true, 1, $false, "Kom ons kyk hoe ParentText werk. ", 2) | Out-Null
$d.Content.Find.Execute("ParentText offers 5 daily lessons to improve your relationship with your child or teen. Once you complete all 5 daily lessons, you will earn a positive parenting certificate! ", $true, $false, $false, $false, $false, $true, 1, $false, "ParentText bied 5 daaglikse lesse aan om jou verhouding met jou kind of tiener te verbeter. Sodra jy al 5 daaglikse lesse voltooi het, sal jy 'n positiewe ouerskapsertifikaat ontvang! ", 2) | Out-Null
$d.Content.Find.Execute("Improve My Relationship with My Child or Teen", $true, $false, $false, $false, $false, $true, 1, $false, "Verbeter My Verhouding met My Kind of Tiener", 2) | Out-Null
$d.Content.Find.Execute("Spending One-on-one Time with My Child or Teen ", $true, $false, $false, $false, $false, $true, 1, $false, "Spandeer Een-tot-een Tyd saam met My Kind of Tiener ", 2) | Out-Null
$d.Content.Find.Execute("Giving Praise ", $true, $false, $false, $false, $false, $true, 1, $false, "Om te Prys ", 2) | Out-Null
$d.Content.Find.Execute("Creating a Routine for One-on-one Time ", $true, $false, $false, $false, $false, $true, 1, $false, "Skep 'n Roetine vir Een-tot-een Tyd ", 2) | Out-Null
$d.Content.Find.Execute("Noticing Feelings During One-on-one Time ", $true, $false, $false, $false, $false, $true, 1, $false, "Let op na Gevoelens Tydens Een-tot-een Tyd ", 2) | Out-Null
$d.Content.Find.Execute("Keeping Calm When We Are Stressed ", $true, $false, $false, $false, $false, $true, 1, $false, "Om Kalm te Bly Wanneer Ons Gestres Is ", 2) | Out-Null
$d.Content.Find.Execute("Now, let’s see what a lesson in ParentText looks like. ", $true, $false, $false, $false, $false, $true, 1, $false, "Nou, kom ons kyk hoe 'n les in ParentText lyk. ", 2) | Out-Null
$d.Content.Find.Execute("You'll receive a daily notification to remind you to complete your lesson. And if you miss it, it is also okay! You can always return to ParentText anytime to catch up on your lesson.", $true, $false, $false, $false, $false, $true, 1, $false, "Jy sal 'n daaglikse kennisgewing ontvang om jou te herinner om jou les te voltooi. En as jy dit mis, is dit ook reg! Jy kan enige tyd na ParentText terugkeer om jou les in te haal.", 2) | Out-Null
$d.Content.Find.Execute("Each lesson is a mix of quizzes, comics, tips, and a fun activity to try at home with your child or family.", $true, $false, $false, $false, $false, $true, 1, $false, "Elke les is 'n mengsel van vasvrae, strokiesprente, wenke en 'n lekker aktiwiteit om by die huis saam met jou kind of gesin te probeer.", 2) | Out-Null
$d.Content.Find.Execute("If you are ever stuck or need help, type MENU or HELP at the end of your lessons to get more support. ", $true, $false, $false, $false, $false, $true, 1, $false, "As jy ooit vashak of hulp nodig het, tik KIESLYS of HELP aan die einde van jou lesse om meer ondersteuning te kry. ", 2) | Out-Null
$d.Content.Find.Execute("When you type HELP anytime, you can get information about resources in your community to address family violence, sexual violence, mental health, or other emergencies. ", $true, $false, $false, $false, $false, $true, 1, $false, "Enige tyd wanneer jy HELP tik, kry jy inligting oor hulpbronne in jou gemeenskap om gesinsgeweld, seksuele geweld, geestesgesondheid of ander noodgevalle aan te spreek. ", 2) | Out-Null
$d.Content.Find.Execute("Your information here is safe: Nothing will be shared without your permission and will not be sold for profit. The messages you send are encrypted and locked in a secure server. ", $true, $false, $false, $false, $false, $true, 1, $false, "Jou inligting is veilig hier: Niks sal gedeel word sonder jou toestemming nie en sal nie vir wins verkoop word nie. Die boodskappe wat jy stuur, is geïnkripteer en gesluit in 'n veilige bediener. ", 2) | Out-Null
$d.Content.Find.Execute("Remember, anyone with access to your unlocked phone can view your messages. So, if you send sensitive information and are worried, delete the messages from your phone. ", $true, $false, $false, $false, $false, $true, 1, $false, "Onthou, enige iemand met toegang tot jou ongeslote foon, kan na jou boodskappe kyk. So, as jy sensitiewe inligting stuur en is bekommerd, verwyder die boodskappe van jou foon. ", 2) | Out-Null
$d.Content.Find.Execute("Being here shows how much you care about providing the best support for your child. ", $true, $false, $false, $false, $false, $true, 1, $false, "Deur hier te wees, wys hoeveel jy omgee om die beste ondersteuning vir jou kind te gee. ", 2) | Out-Null
$d.Content.Find.Execute("It is what you do with your child that will really make a difference. ", $true, $false, $false, $false, $false, $true, 1, $false, "Dit is wat jy doen saam met jou kind wat regtig 'n verskil sal maak. ", 2) | Out-Null
$d.Content.Find.Execute("ParentText will provide tips through lessons to help you with your relationship with your child. It is up to you to put these tips into practice!", $true, $false, $false, $false, $false, $true, 1, $false, "ParentText sal deur lesse wenke gee om jou te help met jou verhouding met jou kind. Jy moet dan net hierdie wenke toepas!", 2) | Out-Null
$d.Content.Find.Execute("Thank you so much for listening! You can access this video at any time via MENU. We hope you enjoy your ParentText journey and make the most out of it! ", $true, $false, $false, $false, $false, $true, 1, $false, "Baie dankie dat jy geluister het! Jy kan enige tyd toegang kry tot hierdie video via KIESLYS. Ons hoop jy geniet jou ParentText reis en maak die meeste daarvan! ", 2) | Out-Null
$d.Content.Find.Execute("Hi! How are you feeling right now? Do you have 30 seconds?", $true, $false, $false, $false, $false, $true, 1, $false, "Hallo! Hoe voel jy op hierdie oomblik? Het jy 30 sekondes?", 2) | Out-Null
$d.Content.Find.Execute("Before you get started in the ParentText programme, let's take a quick pause together.", $true, $false, $false, $false, $false, $true, 1, $false, "Voor jy met die ParentText-program begin, kom ons vat saam 'n rustige oomblik.", 2) | Out-Null
$d.Content.Find.Execute("Try to Take a Pause whenever you feel angry, overwhelmed, stressed, or worried.", $true, $false, $false, $false, $false, $true, 1, $false, "Vat 'n Blaaskans wanneer jy kwaad, oorweldig, gespanne of bekommerd voel.", 2) | Out-Null
$d.Content.Find.Execute("You can also Take a Pause with your child or teen!", $true, $false, $false, $false, $false, $true, 1, $false, "Vat 'n Blaaskans saam met jou kind of tiener!", 2) | Out-Null
$d.Content.Find.Execute("Take a pause with your child or teen!", $true, $false, $false, $false, $false, $true, 1, $false, "Vat 'n Blaaskans saam met jou kind of tiener!", 2) | Out-Null
$d.Content.Find.Execute("Sit down somewhere comfortable and close your eyes.", $true, $false, $false, $false, $false, $true, 1, $false, "Sit iewers gemaklik en maak jou oë toe.", 2) | Out-Null
$d.Content.Find.Execute("Take a deeeeeeeep breath.", $true, $false, $false, $false, $false, $true, 1, $false, "Vat 'n dieeeeeep asem.", 2) | Out-Null
$d.Content.Find.Execute("Feel the air moving in, and out, of your body.", $true, $false, $false, $false, $false, $true, 1, $false, "Voel die lug wat in en uit jou liggaam beweeg.", 2) | Out-Null
$d.Content.Find.Execute("and out;", $true, $false, $false, $false, $false, $true, 1, $false, "en uit;", 2) | Out-Null
$d.Content.Find.Execute("Notice how your body feels while you breathe.", $true, $false, $false, $false, $false, $true, 1, $false, "Let op na hoe jou liggaam voel terwyl jy asemhaal.", 2) | Out-Null
$d.Content.Find.Execute("Notice where you feel tension in your body.", $true, $false, $false, $false, $false, $true, 1, $false, "Let op na waar jy spanning in jou liggaam voel.", 2) | Out-Null
$d.Content.Find.Execute("Try to let it relax.", $true, $false, $false, $false, $false, $true, 1, $false, "Probeer om dit te ontspan.", 2) | Out-Null
$d.Content.Find.Execute("When you are ready, open your eyes again.", $true, $false, $false, $false, $false, $true, 1, $false, "Wanneer jy gereed is, maak weer jou oë oop.", 2) | Out-Null
$d.Content.Find.Execute("Now, notice if you are feeling any differently than", $true, $false, $false, $false, $false, $true, 1, $false, "Nou, let op of jy anders voel as", 2) | Out-Null
$d.Content.Find.Execute("when you started this activity.", $true, $false, $false, $false, $false, $true, 1, $false, "toe jy hierdie aktiwiteit begin het.", 2) | Out-Null
$d.Content.Find.Execute("Even a few deep breaths, or connecting with the ground beneath you, can make a difference.", $true, $false, $false, $false, $false, $true, 1, $false, "Selfs 'n paar diep asems, of net om die grond onder jou te voel, kan 'n verskil maak.", 2) | Out-Null

# Standalone "Take a Pause" paragraph (avoid the "{Take a Pause}" section header)
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Take a Pause") {
        $p.Range.Find.Execute("Take a Pause", $true, $false, $false, $false, $false, $true, 1, $false, "Vat 'n Blaaskans", 2) | Out-Null
    }
}
